# Natmi following Dr Hou advice
# Update the LR-pair stats sheet: a new sending cluster ("Neutro") is added,
# and each sending cluster (ECs, FAPs, Neutro, sCs) now has rows for both
# target clusters (FAPs, sCs) instead of just one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: A..T
#   A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#   E Ligand-expressing cells, F Ligand detection rate,
#   G Ligand average expression value, H Ligand total expression value,
#   I Ligand derived specificity of average expression value,
#   J Ligand derived specificity of total expression value,
#   K Receptor-expressing cells, L Receptor detection rate,
#   M Receptor average expression value, N Receptor total expression value,
#   O Receptor derived specificity of average expression value,
#   P Receptor derived specificity of total expression value,
#   Q Edge average expression weight, R Edge total expression weight,
#   S Edge average expression derived specificity,
#   T Edge total expression derived specificity

$rows = @(
    @{ Row=2;  A="ECs";    B="Timp3"; C="Agtr2"; D="FAPs"; E=3; F=1;                  G=114.380483;        H=343.141449;        I=0.6536977130169467;  J=0.6536977130169466;  K=3; L=1;                  M=0.9721403333333333; N=2.916421; O=0.4203533802573702; P=0.4203533802573701; Q=111.1938808704477;  R=1000.744927834029;  S=0.2747840433331858;   T=0.2747840433331857 }
    @{ Row=3;  A="ECs";    B="Timp3"; C="Agtr2"; D="sCs";  E=3; F=1;                  G=114.380483;        H=343.141449;        I=0.6536977130169467;  J=0.6536977130169466;  K=2; L=0.6666666666666666; M=1.340533666666667;  N=4.021601; O=0.5796466197426299; P=0.5796466197426299; Q=153.3308882710944;  R=1379.977994439849;  S=0.3789136696837609;   T=0.3789136696837608 }
    @{ Row=4;  A="FAPs";   B="Timp3"; C="Agtr2"; D="FAPs"; E=3; F=1;                  G=15.70856733333333; H=47.125702;          I=0.08977628238003432; J=0.08977628238003432; K=3; L=1;                  M=0.9721403333333333; N=2.916421; O=0.4203533802573702; P=0.4203533802573701; Q=15.27093188361578;  R=137.438386952542;   S=0.03773776376538761;  T=0.03773776376538761 }
    @{ Row=5;  A="FAPs";   B="Timp3"; C="Agtr2"; D="sCs";  E=3; F=1;                  G=15.70856733333333; H=47.125702;          I=0.08977628238003432; J=0.08977628238003432; K=2; L=0.6666666666666666; M=1.340533666666667;  N=4.021601; O=0.5796466197426299; P=0.5796466197426299; Q=21.05786336543356;  R=189.520770288902;   S=0.05203851861464671;  T=0.05203851861464671 }
    @{ Row=6;  A="Neutro"; B="Timp3"; C="Agtr2"; D="FAPs"; E=2; F=0.6666666666666666; G=0.1510706666666667; H=0.453212;          I=0.0008633863637727903; J=0.0008633863637727903; K=3; L=1;                  M=0.9721403333333333; N=2.916421; O=0.4203533802573702; P=0.4203533802573701; Q=0.1468618882502222; R=1.321756994252;     S=0.0003629273764800119; T=0.0003629273764800118 }
    @{ Row=7;  A="Neutro"; B="Timp3"; C="Agtr2"; D="sCs";  E=2; F=0.6666666666666666; G=0.1510706666666667; H=0.453212;          I=0.0008633863637727903; J=0.0008633863637727903; K=2; L=0.6666666666666666; M=1.340533666666667;  N=4.021601; O=0.5796466197426299; P=0.5796466197426299; Q=0.2025153147124444; R=1.822637832412;     S=0.0005004589872927785; T=0.0005004589872927785 }
    @{ Row=8;  A="sCs";    B="Timp3"; C="Agtr2"; D="FAPs"; E=3; F=1;                  G=44.73445933333333; H=134.203378;         I=0.2556626182392462;  J=0.2556626182392462;  K=3; L=1;                  M=0.9721403333333333; N=2.916421; O=0.4203533802573702; P=0.4203533802573701; Q=43.4881722077931;   R=391.3935498701379;  S=0.1074686457823167;   T=0.1074686457823167 }
    @{ Row=9;  A="sCs";    B="Timp3"; C="Agtr2"; D="sCs";  E=3; F=1;                  G=44.73445933333333; H=134.203378;         I=0.2556626182392462;  J=0.2556626182392462;  K=2; L=0.6666666666666666; M=1.340533666666667;  N=4.021601; O=0.5796466197426299; P=0.5796466197426299; Q=59.96804879646422;  R=539.712439168178;   S=0.1481939724569295;   T=0.1481939724569295 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
